# Generate Report for Handoff
# Replace the old localization file GUID/hash stamp with the new one, and
# bump the "Latest Handoff" timestamps to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

$oldGuid = "9f27db05-905b-4c85-9940-869ca5bf0b6c"
$newGuid = "14d30e09-5902-4338-a1a8-2cf0d70ead89"

$oldHash = "c4d66aba43c2eb95dadbf4fbfe998014791a8be1"
$newHash = "0632078ed0787b45800ec8667fae4eae79f47fed"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-13 15:11:16"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-13 15:11:07"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-13 15:11:16"
